# Align WHO Influenza observation codes with current WHO guidance.
# Adds four new coded-observation rows (1023-1026) to the "Conditions"
# sheet of the WHO Coded Observations workbook, and brings the existing
# "Healthcare worker / Older adult / Immunocompromised" rows (10-12) into
# the same formatting as the rest of the table (rows 1-9) rather than
# their previous one-off style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conditions")

# --- 1. Re-format rows 10-12 (Healthcare worker / Older adult 60+ /
#        Immunocompromised) so they match the table's normal row style
#        instead of the old one-off style. ---
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H12").PasteSpecial(-4122)

# --- 2. Append the four new WHO influenza risk-indication rows. ---
# Row 13: Pregnant
$ws.Range("A13").Value = "1023"
$ws.Range("B13").Value = "Pregnant"
$ws.Range("C13").Value = "Patient is pregnant"
$ws.Range("D13").Value = "n/a"
$ws.Range("E13").Value = "n/a"
$ws.Range("F13").Value = "77386006"
$ws.Range("G13").Value = "n/a"
$ws.Range("H13").Value = "n/a"

# Row 14: Chronic medical condition
$ws.Range("A14").Value = "1024"
$ws.Range("B14").Value = "Chronic medical condition"
$ws.Range("C14").Value = "Patient has a chronic medical condition increasing influenza risk"
$ws.Range("D14").Value = "n/a"
$ws.Range("E14").Value = "n/a"
$ws.Range("F14").Value = "27624003"
$ws.Range("G14").Value = "n/a"
$ws.Range("H14").Value = "n/a"

# Row 15: WHO influenza priority - young child
$ws.Range("A15").Value = "1025"
$ws.Range("B15").Value = "WHO influenza priority - young child"
$ws.Range("C15").Value = "Child aged 6 months to 5 years (WHO influenza priority group)"
$ws.Range("D15").Value = "n/a"
$ws.Range("E15").Value = "n/a"
$ws.Range("F15").Value = "410601007"
$ws.Range("G15").Value = "n/a"
$ws.Range("H15").Value = "n/a"

# Row 16: WHO influenza priority - older adult
$ws.Range("A16").Value = "1026"
$ws.Range("B16").Value = "WHO influenza priority - older adult"
$ws.Range("C16").Value = "Older adult aged 65 years or older (WHO influenza priority group)"
$ws.Range("D16").Value = "n/a"
$ws.Range("E16").Value = "n/a"
$ws.Range("F16").Value = "105436006"
$ws.Range("G16").Value = "n/a"
$ws.Range("H16").Value = "n/a"

# Give the new rows the same (pre-reformat) look the observation rows
# used to have, matched here against row 10's now-normalized style so
# the new block reads consistently with the rest of the table.
$ws.Range("A10:H10").Copy()
$ws.Range("A13:H16").PasteSpecial(-4122)

# Re-apply the text values (PasteSpecial(xlPasteFormats) only touches
# formatting, but make sure the number format stays general text so the
# codes such as "1023" are not reinterpreted as numbers).
$ws.Range("A13").Value = "1023"
$ws.Range("B13").Value = "Pregnant"
$ws.Range("C13").Value = "Patient is pregnant"
$ws.Range("D13").Value = "n/a"
$ws.Range("E13").Value = "n/a"
$ws.Range("F13").Value = "77386006"
$ws.Range("G13").Value = "n/a"
$ws.Range("H13").Value = "n/a"

$ws.Range("A14").Value = "1024"
$ws.Range("B14").Value = "Chronic medical condition"
$ws.Range("C14").Value = "Patient has a chronic medical condition increasing influenza risk"
$ws.Range("D14").Value = "n/a"
$ws.Range("E14").Value = "n/a"
$ws.Range("F14").Value = "27624003"
$ws.Range("G14").Value = "n/a"
$ws.Range("H14").Value = "n/a"

$ws.Range("A15").Value = "1025"
$ws.Range("B15").Value = "WHO influenza priority - young child"
$ws.Range("C15").Value = "Child aged 6 months to 5 years (WHO influenza priority group)"
$ws.Range("D15").Value = "n/a"
$ws.Range("E15").Value = "n/a"
$ws.Range("F15").Value = "410601007"
$ws.Range("G15").Value = "n/a"
$ws.Range("H15").Value = "n/a"

$ws.Range("A16").Value = "1026"
$ws.Range("B16").Value = "WHO influenza priority - older adult"
$ws.Range("C16").Value = "Older adult aged 65 years or older (WHO influenza priority group)"
$ws.Range("D16").Value = "n/a"
$ws.Range("E16").Value = "n/a"
$ws.Range("F16").Value = "105436006"
$ws.Range("G16").Value = "n/a"
$ws.Range("H16").Value = "n/a"
